$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003455333333333333
$ws.Range("H2").Value = 0.010366
$ws.Range("I2").Value = 0.000270121469710956
$ws.Range("J2").Value = 0.000270121469710956
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.022792
$ws.Range("N2").Value = 0.06837600000000001
$ws.Range("O2").Value = 0.001916327914826657
$ws.Range("P2").Value = 0.001916327914826657
$ws.Range("Q2").Value = 0.00007875395733333334
$ws.Range("R2").Value = 0.0007087856160000001
$ws.Range("S2").Value = 0.0000005176413128011083
$ws.Range("T2").Value = 0.0000005176413128011084

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003455333333333333
$ws.Range("H3").Value = 0.010366
$ws.Range("I3").Value = 0.000270121469710956
$ws.Range("J3").Value = 0.000270121469710956
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.402094666666667
$ws.Range("N3").Value = 13.206284
$ws.Range("O3").Value = 0.3701235913233977
$ws.Range("P3").Value = 0.3701235913233977
$ws.Range("Q3").Value = 0.01521070443822222
$ws.Range("R3").Value = 0.136896339944
$ws.Range("S3").Value = 0.00009997832846297342
$ws.Range("T3").Value = 0.00009997832846297344

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.003455333333333333
$ws.Range("H4").Value = 0.010366
$ws.Range("I4").Value = 0.000270121469710956
$ws.Range("J4").Value = 0.000270121469710956
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.468693666666667
$ws.Range("N4").Value = 22.406081
$ws.Range("O4").Value = 0.6279600807617757
$ws.Range("P4").Value = 0.6279600807617757
$ws.Range("Q4").Value = 0.02580682618288889
$ws.Range("R4").Value = 0.232261435646
$ws.Range("S4").Value = 0.0001696254999351815
$ws.Range("T4").Value = 0.0001696254999351815

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.502271
$ws.Range("H5").Value = 34.50681299999999
$ws.Range("I5").Value = 0.8991926531546518
$ws.Range("J5").Value = 0.8991926531546519
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.022792
$ws.Range("N5").Value = 0.06837600000000001
$ws.Range("O5").Value = 0.001916327914826657
$ws.Range("P5").Value = 0.001916327914826657
$ws.Range("Q5").Value = 0.262159760632
$ws.Range("R5").Value = 2.359437845688
$ws.Range("S5").Value = 0.001723147982047303
$ws.Range("T5").Value = 0.001723147982047304

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.502271
$ws.Range("H6").Value = 34.50681299999999
$ws.Range("I6").Value = 0.8991926531546518
$ws.Range("J6").Value = 0.8991926531546519
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.402094666666667
$ws.Range("N6").Value = 13.206284
$ws.Range("O6").Value = 0.3701235913233977
$ws.Range("P6").Value = 0.3701235913233977
$ws.Range("Q6").Value = 50.63408582365466
$ws.Range("R6").Value = 455.706772412892
$ws.Range("S6").Value = 0.3328124140772141
$ws.Range("T6").Value = 0.3328124140772141

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.502271
$ws.Range("H7").Value = 34.50681299999999
$ws.Range("I7").Value = 0.8991926531546518
$ws.Range("J7").Value = 0.8991926531546519
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.468693666666667
$ws.Range("N7").Value = 22.406081
$ws.Range("O7").Value = 0.6279600807617757
$ws.Range("P7").Value = 0.6279600807617757
$ws.Range("Q7").Value = 85.90693856998365
$ws.Range("R7").Value = 773.1624471298529
$ws.Range("S7").Value = 0.5646570910953905
$ws.Range("T7").Value = 0.5646570910953905

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.286049666666667
$ws.Range("H8").Value = 3.858149
$ws.Range("I8").Value = 0.1005372253756372
$ws.Range("J8").Value = 0.1005372253756372
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.022792
$ws.Range("N8").Value = 0.06837600000000001
$ws.Range("O8").Value = 0.001916327914826657
$ws.Range("P8").Value = 0.001916327914826657
$ws.Range("Q8").Value = 0.02931164400266667
$ws.Range("R8").Value = 0.263804796024
$ws.Range("S8").Value = 0.0001926622914665525
$ws.Range("T8").Value = 0.0001926622914665525

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.286049666666667
$ws.Range("H9").Value = 3.858149
$ws.Range("I9").Value = 0.1005372253756372
$ws.Range("J9").Value = 0.1005372253756372
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.402094666666667
$ws.Range("N9").Value = 13.206284
$ws.Range("O9").Value = 0.3701235913233977
$ws.Range("P9").Value = 0.3701235913233977
$ws.Range("Q9").Value = 5.661312378701778
$ws.Range("R9").Value = 50.951811408316
$ws.Range("S9").Value = 0.03721119891772067
$ws.Range("T9").Value = 0.03721119891772068

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.286049666666667
$ws.Range("H10").Value = 3.858149
$ws.Range("I10").Value = 0.1005372253756372
$ws.Range("J10").Value = 0.1005372253756372
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.468693666666667
$ws.Range("N10").Value = 22.406081
$ws.Range("O10").Value = 0.6279600807617757
$ws.Range("P10").Value = 0.6279600807617757
$ws.Range("Q10").Value = 9.605111000452112
$ws.Range("R10").Value = 86.44599900406901
$ws.Range("S10").Value = 0.06313336416644998
$ws.Range("T10").Value = 0.06313336416644999
